$d = $word.ActiveDocument

# Locate the target paragraph (the last paragraph of "Results of Experiments"),
# found by its distinctive leading text rather than a hard-coded index.
$target = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "Overall very few difficulties occurred*") {
        $target = $p
    }
}
if ($target -eq $null) {
    throw "Could not locate target paragraph"
}

$paraRange = $d.Range($target.Range.Start, $target.Range.End)
$insertPos = $paraRange.Start

# Remove the whole paragraph (text, runs, and the trailing _GoBack bookmark
# that lives at its end) -- it will be rebuilt (split across several new
# paragraphs) from scratch below.
$paraRange.Delete()

$fragment = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
      <w:pPr>
        <w:rPr>
          <w:rFonts w:eastAsia="Times New Roman" w:cs="Times New Roman"/>
          <w:lang w:eastAsia="en-US"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:t>Overall very few difficulties occurred in the actual development of the system. We hit a few road blocks in design, but we worked that out at the following meeting as noted above.</w:t>
      </w:r>
      <w:r>
        <w:t xml:space="preserve"> </w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:eastAsia="Times New Roman" w:cs="Times New Roman"/>
          <w:lang w:eastAsia="en-US"/>
        </w:rPr>
        <w:t>On the technical side, many of our roadblocks involved timing, deadlock prevention, and synchroniz</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:eastAsia="Times New Roman" w:cs="Times New Roman"/>
          <w:lang w:eastAsia="en-US"/>
        </w:rPr>
        <w:t xml:space="preserve">ing the threads during the day. </w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:eastAsia="Times New Roman" w:cs="Times New Roman"/>
          <w:lang w:eastAsia="en-US"/>
        </w:rPr>
        <w:t>The timing was fixed by using a starting countdown latch on all threads and by changing the interval tha</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:eastAsia="Times New Roman" w:cs="Times New Roman"/>
          <w:lang w:eastAsia="en-US"/>
        </w:rPr>
        <w:t>t threads wait when doing work.</w:t>
      </w:r>
    </w:p>
    <w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
      <w:pPr>
        <w:rPr>
          <w:rFonts w:eastAsia="Times New Roman" w:cs="Times New Roman"/>
          <w:lang w:eastAsia="en-US"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:rFonts w:eastAsia="Times New Roman" w:cs="Times New Roman"/>
          <w:lang w:eastAsia="en-US"/>
        </w:rPr>
        <w:t>The conditions for deadlo</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:eastAsia="Times New Roman" w:cs="Times New Roman"/>
          <w:lang w:eastAsia="en-US"/>
        </w:rPr>
        <w:t>ck to occur existed in our code</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:eastAsia="Times New Roman" w:cs="Times New Roman"/>
          <w:lang w:eastAsia="en-US"/>
        </w:rPr>
        <w:t xml:space="preserve"> in a synchronized method</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:eastAsia="Times New Roman" w:cs="Times New Roman"/>
          <w:lang w:eastAsia="en-US"/>
        </w:rPr>
        <w:t xml:space="preserve"> near the end of the day. All employees countdown and await the 4:00 meeting countdow</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:eastAsia="Times New Roman" w:cs="Times New Roman"/>
          <w:lang w:eastAsia="en-US"/>
        </w:rPr>
        <w:t>n latch</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:eastAsia="Times New Roman" w:cs="Times New Roman"/>
          <w:lang w:eastAsia="en-US"/>
        </w:rPr>
        <w:t>.</w:t>
      </w:r>
      <w:bookmarkStart w:id="0" w:name="_GoBack"/>
      <w:bookmarkEnd w:id="0"/>
      <w:r>
        <w:rPr>
          <w:rFonts w:eastAsia="Times New Roman" w:cs="Times New Roman"/>
          <w:lang w:eastAsia="en-US"/>
        </w:rPr>
        <w:t xml:space="preserve"> Deadlock could occur if the Project Manager entered this method, holding its own lock, while an employee tried to ask a question. </w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:eastAsia="Times New Roman" w:cs="Times New Roman"/>
          <w:lang w:eastAsia="en-US"/>
        </w:rPr>
        <w:lastRenderedPageBreak/>
        <w:t>The employee would wait for the Project Manager to release its lock, but that would never occur until the employee's question was answered. We solved this by removing the synch</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:eastAsia="Times New Roman" w:cs="Times New Roman"/>
          <w:lang w:eastAsia="en-US"/>
        </w:rPr>
        <w:t>ronization</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:eastAsia="Times New Roman" w:cs="Times New Roman"/>
          <w:lang w:eastAsia="en-US"/>
        </w:rPr>
        <w:t>. </w:t>
      </w:r>
    </w:p>
    <w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
      <w:pPr>
        <w:rPr>
          <w:rFonts w:eastAsia="Times New Roman" w:cs="Times New Roman"/>
          <w:lang w:eastAsia="en-US"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:rFonts w:eastAsia="Times New Roman" w:cs="Times New Roman"/>
          <w:lang w:eastAsia="en-US"/>
        </w:rPr>
        <w:t>Synchronizing the threads during the day became complicated as more countdown latches were needed at different parts of the da</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:eastAsia="Times New Roman" w:cs="Times New Roman"/>
          <w:lang w:eastAsia="en-US"/>
        </w:rPr>
        <w:t>y. We considered using barriers;</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:eastAsia="Times New Roman" w:cs="Times New Roman"/>
          <w:lang w:eastAsia="en-US"/>
        </w:rPr>
        <w:t xml:space="preserve"> however, many of the activit</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:eastAsia="Times New Roman" w:cs="Times New Roman"/>
          <w:lang w:eastAsia="en-US"/>
        </w:rPr>
        <w:t xml:space="preserve">ies in the day happen between </w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:eastAsia="Times New Roman" w:cs="Times New Roman"/>
          <w:lang w:eastAsia="en-US"/>
        </w:rPr>
        <w:t xml:space="preserve">distinct </w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:eastAsia="Times New Roman" w:cs="Times New Roman"/>
          <w:lang w:eastAsia="en-US"/>
        </w:rPr>
        <w:t>groups</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:eastAsia="Times New Roman" w:cs="Times New Roman"/>
          <w:lang w:eastAsia="en-US"/>
        </w:rPr>
        <w:t xml:space="preserve"> of employees once a day. This increased coupling in our code, as employees need to know all the latches they will use in a day. </w:t>
      </w:r>
    </w:p>
    <w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
      <w:pPr>
        <w:rPr>
          <w:rFonts w:eastAsia="Times New Roman" w:cs="Times New Roman"/>
          <w:lang w:eastAsia="en-US"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:rFonts w:eastAsia="Times New Roman" w:cs="Times New Roman"/>
          <w:lang w:eastAsia="en-US"/>
        </w:rPr>
        <w:t xml:space="preserve">We did not have roadblocks with providing access to things like the conference room, Team Leads, and the Project Manager. This is because we carefully designed the project so as to easily allow synchronized access to </w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:eastAsia="Times New Roman" w:cs="Times New Roman"/>
          <w:lang w:eastAsia="en-US"/>
        </w:rPr>
        <w:t>parts</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:eastAsia="Times New Roman" w:cs="Times New Roman"/>
          <w:lang w:eastAsia="en-US"/>
        </w:rPr>
        <w:t xml:space="preserve"> requiring it. </w:t>
      </w:r>
    </w:p>
    <w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"/>
'@

$insertRange = $d.Range($insertPos, $insertPos)
$insertRange.InsertXML($fragment)

Write-Output ("Paragraph count after edit: " + $d.Paragraphs.Count)
